$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "285.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-5.51%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.21%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.953"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.55%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07254"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-7.18%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.777"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-17.41%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.650"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.95%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9105"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.06%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1634"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.99%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07492"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.78%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08163"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.03%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.73%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09999"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.21%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001504"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.13%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005747"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.87%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.007498"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,116.77%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.469"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.20%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.722"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.70%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.121"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.41%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3259"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.03%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1294"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.58%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.385"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.41%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04511"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.85%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001241"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.03%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003991"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.89%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001264"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "1.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-9.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04324"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-9.43%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007448"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.47%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1308"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.46%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002166"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.74%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01073"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.79%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006090"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.31%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000759"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.12%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.896"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "132.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003035"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-14.75%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002124"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.12%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002023"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.12%"
